$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.487.70"
$ws.Range("E2").Value = "  +6.15%  "
$ws.Range("D3").Value = "2.046.19"
$ws.Range("E3").Value = "  +3.34%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.36"
$ws.Range("E5").Value = "  +5.19%  "
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.23"
$ws.Range("E7").Value = "  +18.48%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  +5.96%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.49"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +4.47%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.911"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.19"
$ws.Range("E14").Value = "  +7.13%  "
$ws.Range("D15").Value = "2.348.63"
$ws.Range("E15").Value = "  +3.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.63"
$ws.Range("E16").Value = "  +7.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.74"
$ws.Range("E17").Value = "  +21.59%  "
$ws.Range("D18").Value = "2.055.78"
$ws.Range("E18").Value = "  +3.82%  "
$ws.Range("D19").Value = "37.401.67"
$ws.Range("E19").Value = "  +6.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.55"
$ws.Range("E20").Value = "  +5.42%  "
$ws.Range("D21").Value = "0.0₃0880"
$ws.Range("E21").Value = "  +5.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.36"
$ws.Range("E22").Value = "  +7.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.58"
$ws.Range("E23").Value = "  +2.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.73"
$ws.Range("E24").Value = "  +22.19%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("E26").Value = "  +4.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.60"
$ws.Range("E27").Value = "  +6.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.53"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.26"
$ws.Range("E30").Value = "  +11.11%  "
$ws.Range("E31").Value = "  +3.22%  "
$ws.Range("E32").Value = "  +7.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("E33").Value = "  +22.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.73"
$ws.Range("E34").Value = "  +12.02%  "
$ws.Range("E35").Value = "  +5.63%  "
$ws.Range("E36").Value = "  +9.40%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.04"
$ws.Range("E38").Value = "  +25.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.81"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("E40").Value = "  +17.75%  "
$ws.Range("E41").Value = "  +5.37%  "
$ws.Range("E42").Value = "  +5.56%  "
$ws.Range("E43").Value = "  +6.10%  "
$ws.Range("E44").Value = "  +6.80%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.72"
$ws.Range("E45").Value = "  +21.49%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.10"
$ws.Range("E46").Value = "  +9.52%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.04"
$ws.Range("E47").Value = "  +11.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.51"
$ws.Range("E48").Value = "  +6.62%  "
$ws.Range("D49").Value = "1.429.56"
$ws.Range("E49").Value = "  +5.63%  "
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.57"
$ws.Range("E51").Value = "  +5.33%  "
